# Week8_Recap.pptx – fix the off-by-one in the array-copy loop shown on the
# "String" slide: the loop condition comment changes from
#   i <= 7
# to
#   i < 7
# (the code sample was copying one element too many).
#
# The text lives in slide 11 (sldId 647), shape Id=13 ("TextBox 12",
# creationId {2B89EFD4-1CE0-424E-8F8A-D7F3C08BDEC9}), inside the run
# that reads " <= ". A visually similar " <= n; " snippet also exists on
# another slide, so we search defensively and only touch the shape that
# matches both the expected Id/Name *and* contains the exact substring.

$p = $ppt.ActivePresentation

$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.Id -eq 13 -and $shape.Name -eq "TextBox 12" -and $shape.HasTextFrame) {
            $candidateText = $shape.TextFrame.TextRange.Text
            if ($candidateText.IndexOf(" <= ") -ge 0) {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$matchIndex = $fullText.IndexOf(" <= ")

# Characters() is 1-based; replace just the " <= " run text with " < ".
$run = $tr.Characters($matchIndex + 1, 4)
$run.Text = " < "
